# Add a new bash_lib entry (row 63) documenting the BSD/UNIX/GNU parameter
# styles, as described in the commit "Add bash entry on command style:
# BSD/UNIX/GNU".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A repeats the "linux" category tag used by the surrounding rows,
# column B is the short title, column C holds the (wrapped) body text.
$ws.Range("A63").Value = "linux"
$ws.Range("B63").Value = "parameter style"
$ws.Range("C63").Value = "ps aux    //this is BSD style`nps -elf    //this is UNIX Style`ngrep --color    //this is GNU style"

# Column C uses the same wrapped-text style as every other row in the table.
$ws.Range("C63").WrapText = $true

# Row height matches the new wrapped content (3 lines).
$ws.Rows.Item(63).RowHeight = 47.25

# Move the view down to the newly added row and select the next empty cell,
# mirroring where Excel would land the cursor after typing the new row.
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C64").Select()
